$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.011.75"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.557.58"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.33"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.487"
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.69"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0590"
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.780.46"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.558.67"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.018.69"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.99"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "216.04"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.28"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.23"
$ws.Range("E23").Value = "  +2.61%  "
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.78"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.93"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.402.71"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.00"
$ws.Range("E34").Value = "  +3.30%  "
$ws.Range("E35").Value = "  +3.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.962"
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.523"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.811"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.992"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("E43").Value = "  +3.32%  "
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.94"
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.694.69"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.26"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0961"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("E51").Value = "  +0.40%  "
